# Edit: add dataset description text to the "Source/Description/Volume/Variety" textbox
# on slide 4, remove the stray "Velocity/Variability" paragraph, and resize the textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$cr = [char]13

# Build the full replacement text (8 paragraphs, separated by carriage returns).
# Paragraph layout:
#   1. "Source: " + "Rotten Tomatoes (Kaggle) and IMDB (IMDB)"
#   2. (empty)
#   3. "Description: " + long description
#   4. (empty)
#   5. "Volume: " + "~174MB"
#   6. (empty)
#   7. "Variety" + ": IMDB's dataset used to retrieve information about the films and Rotten Tomatoes' one used to obtain the reviews about the movies"
#   8. (empty, final)

$rsquo = [char]8217

$sourceLabel = "Source: "
$sourceValue = "Rotten Tomatoes (Kaggle) and IMDB (IMDB)"

$descLabel = "Description: "
$descValue = "The dataset contains basic information about the movies like the title, the cast and crew with roles, the year of release, the runtime and the production houses. In addition, the dataset contains the reviews and scores for the movies divided between top critics and audience.  "

$volumeLabel = "Volume: "
$volumeValue = "~174MB"

$varietyLabel = "Variety"
$varietyValue = ": IMDB" + $rsquo + "s dataset used to retrieve information about the films and Rotten Tomatoes" + $rsquo + " one used to obtain the reviews about the movies"

$fullText = $sourceLabel + $sourceValue + $cr + $cr + $descLabel + $descValue + $cr + $cr + $volumeLabel + $volumeValue + $cr + $cr + $varietyLabel + $varietyValue + $cr

$tr.Text = $fullText

# ---- Paragraph 1: "Source: " (bold italic) + value (colored, non-italic, custom font) ----
$start = 1
$len = $sourceLabel.Length
$runSourceLabel = $tr.Characters($start, $len)
$runSourceLabel.Font.Bold = 1
$runSourceLabel.Font.Italic = 1

$start = $start + $len
$len = $sourceValue.Length
$runSourceValue = $tr.Characters($start, $len)
$runSourceValue.Font.Bold = 0
$runSourceValue.Font.Italic = 0
$runSourceValue.Font.Name = "zeitung"
$runSourceValue.Font.Color.RGB = 2367776
$runSourceValue.Font.Shadow = 0

# ---- Paragraph 3: "Description: " (bold italic) + value (plain) ----
$start = $start + $len + 2
$len = $descLabel.Length
$runDescLabel = $tr.Characters($start, $len)
$runDescLabel.Font.Bold = 1
$runDescLabel.Font.Italic = 1

$start = $start + $len
$len = $descValue.Length
$runDescValue = $tr.Characters($start, $len)
$runDescValue.Font.Bold = 0
$runDescValue.Font.Italic = 0

# ---- Paragraph 5: "Volume: " (bold italic) + value (plain) ----
$start = $start + $len + 2
$len = $volumeLabel.Length
$runVolumeLabel = $tr.Characters($start, $len)
$runVolumeLabel.Font.Bold = 1
$runVolumeLabel.Font.Italic = 1

$start = $start + $len
$len = $volumeValue.Length
$runVolumeValue = $tr.Characters($start, $len)
$runVolumeValue.Font.Bold = 0
$runVolumeValue.Font.Italic = 0

# ---- Paragraph 7: "Variety" (bold italic) + value (plain) ----
$start = $start + $len + 2
$len = $varietyLabel.Length
$runVarietyLabel = $tr.Characters($start, $len)
$runVarietyLabel.Font.Bold = 1
$runVarietyLabel.Font.Italic = 1

$start = $start + $len
$len = $varietyValue.Length
$runVarietyValue = $tr.Characters($start, $len)
$runVarietyValue.Font.Bold = 0
$runVarietyValue.Font.Italic = 0

# ---- Resize the textbox to fit the new (longer) content ----
# Shape.Top/Left/Width/Height are expressed in points (1 pt = 12700 EMU);
# only the height changes (8349029 x 3785652 EMU box at the same position).
# Height is stored internally as a single-precision float, so a plain
# 3785652/12700 division truncates one EMU short on save; nudge it
# (well within the float32 rounding margin) to land on the exact EMU value.
$sh.Height = 298.082870
